# Apply the "cryptos list" data refresh described by the commit:
#   "Updated cryptos list on Tue Oct 31 04:12:34 UTC 2023 with GitHub Actions"
#
# Each row holds a coins Price (col D) and Volume(1h) change (col E) as plain
# text (the sheet stores these as strings, e.g. "34.313.64", "  -0.11%  ").
#
# Excel auto-converts plain-looking decimals (e.g. "227.32") typed into a cell
# into a real number. To keep such values as text (matching the original files
# inline-string cells) we enter them with a leading apostrophe (forces "quote
# prefixed" text entry) and then reset the cell Style back to "Normal" so no
# stray number-format / quote-prefix styling is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2: Bitcoin
$ws.Range("D2").Value = "34.313.64"
$ws.Range("E2").Value = "  -0.11%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.800.50"
$ws.Range("E3").Value = "  +0.74%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.13%  "

# Row 5: BNB
$ws.Range("D5").Value = "'227.32"   # force text, avoid numeric auto-conversion
$ws.Range("D5").Style = "Normal"        # drop the quote-prefix style Excel just added

# Row 6: XRP
$ws.Range("E6").Value = "  +3.72%  "

# Row 7: USDC
$ws.Range("E7").Value = "  +0.12%  "

# Row 8: Solana
$ws.Range("D8").Value = "'35.79"   # force text, avoid numeric auto-conversion
$ws.Range("D8").Style = "Normal"        # drop the quote-prefix style Excel just added
$ws.Range("E8").Value = "  +8.85%  "

# Row 9: Cardano
$ws.Range("E9").Value = "  +1.90%  "

# Row 10: Dogecoin
$ws.Range("E10").Value = "  +0.42%  "

# Row 11: TRON
$ws.Range("E11").Value = "  +1.90%  "

# Row 12: WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "2.060.83"
$ws.Range("E12").Value = "  +0.78%  "

# Row 13: Chainlink
$ws.Range("D13").Value = "'11.45"   # force text, avoid numeric auto-conversion
$ws.Range("D13").Style = "Normal"        # drop the quote-prefix style Excel just added
$ws.Range("E13").Value = "  +2.24%  "

# Row 14: WrappedEther
$ws.Range("D14").Value = "1.801.47"
$ws.Range("E14").Value = "  +0.64%  "

# Row 15: Polygon
$ws.Range("D15").Value = "'0.643"   # force text, avoid numeric auto-conversion
$ws.Range("D15").Style = "Normal"        # drop the quote-prefix style Excel just added
$ws.Range("E15").Value = "  +1.43%  "

# Row 16: Polkadot
$ws.Range("D16").Value = "'4.49"   # force text, avoid numeric auto-conversion
$ws.Range("D16").Style = "Normal"        # drop the quote-prefix style Excel just added
$ws.Range("E16").Value = "  +4.73%  "

# Row 17: WrappedBTC
$ws.Range("D17").Value = "34.319.74"
$ws.Range("E17").Value = "  -0.04%  "

# Row 18: Litecoin
$ws.Range("D18").Value = "'68.97"   # force text, avoid numeric auto-conversion
$ws.Range("D18").Style = "Normal"        # drop the quote-prefix style Excel just added
$ws.Range("E18").Value = "  +0.86%  "

# Row 19: BitcoinCash
$ws.Range("D19").Value = "'245.24"   # force text, avoid numeric auto-conversion
$ws.Range("D19").Style = "Normal"        # drop the quote-prefix style Excel just added
$ws.Range("E19").Value = "  +0.11%  "

# Row 20: ShibaInu
$ws.Range("E20").Value = "  -0.05%  "

# Row 21: Avalanche
$ws.Range("D21").Value = "'11.44"   # force text, avoid numeric auto-conversion
$ws.Range("D21").Style = "Normal"        # drop the quote-prefix style Excel just added
$ws.Range("E21").Value = "  +1.71%  "

# Row 22: Dai
$ws.Range("E22").Value = "  +0.06%  "

# Row 23: Uniswap
$ws.Range("E23").Value = "  +0.62%  "

# Row 24: Toncoin->Monero (row content swap)
$ws.Range("B24").Value = "Monero"
$ws.Range("C24").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D24").Value = "'171.47"   # force text, avoid numeric auto-conversion
$ws.Range("D24").Style = "Normal"        # drop the quote-prefix style Excel just added
$ws.Range("E24").Value = "  +1.48%  "

# Row 25: Monero->Toncoin (row content swap)
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "'2.13"   # force text, avoid numeric auto-conversion
$ws.Range("D25").Style = "Normal"        # drop the quote-prefix style Excel just added
$ws.Range("E25").Value = "  +3.11%  "

# Row 26: Cosmos
$ws.Range("E26").Value = "  +7.31%  "

# Row 27: EthereumClassic
$ws.Range("D27").Value = "'16.80"   # force text, avoid numeric auto-conversion
$ws.Range("D27").Style = "Normal"        # drop the quote-prefix style Excel just added
$ws.Range("E27").Value = "  +1.70%  "

# Row 28: Stellar
$ws.Range("E28").Value = "  +2.56%  "

# Row 29: BinanceUSD
$ws.Range("E29").Value = "  +0.08%  "

# Row 30: InternetComputer(DFINITY)
$ws.Range("D30").Value = "'4.03"   # force text, avoid numeric auto-conversion
$ws.Range("D30").Style = "Normal"        # drop the quote-prefix style Excel just added
$ws.Range("E30").Value = "  +0.35%  "

# Row 31: Hedera
$ws.Range("D31").Value = "'0.0530"   # force text, avoid numeric auto-conversion
$ws.Range("D31").Style = "Normal"        # drop the quote-prefix style Excel just added
$ws.Range("E31").Value = "  +0.71%  "

# Row 32: Filecoin
$ws.Range("E32").Value = "  +1.21%  "

# Row 33: PancakeSwap
$ws.Range("E33").Value = "  +0.71%  "

# Row 34: LidoDAOToken
$ws.Range("E34").Value = "  +0.31%  "

# Row 35: Maker
$ws.Range("D35").Value = "1.394.70"
$ws.Range("E35").Value = "  -1.20%  "

# Row 36: ImmutableX
$ws.Range("D36").Value = "'0.671"   # force text, avoid numeric auto-conversion
$ws.Range("D36").Style = "Normal"        # drop the quote-prefix style Excel just added
$ws.Range("E36").Value = "  -1.38%  "

# Row 37: RenderToken
$ws.Range("E37").Value = "  -4.92%  "

# Row 38: TrustWalletToken
$ws.Range("E38").Value = "  -0.44%  "

# Row 39: VeChain
$ws.Range("D39").Value = "'0.0189"   # force text, avoid numeric auto-conversion
$ws.Range("D39").Style = "Normal"        # drop the quote-prefix style Excel just added
$ws.Range("E39").Value = "  -0.54%  "

# Row 40: WEMIXToken
$ws.Range("D40").Value = "'1.24"   # force text, avoid numeric auto-conversion
$ws.Range("D40").Style = "Normal"        # drop the quote-prefix style Excel just added
$ws.Range("E40").Value = "  +11.61%  "

# Row 41: ARBITRUM
$ws.Range("D41").Value = "'0.957"   # force text, avoid numeric auto-conversion
$ws.Range("D41").Style = "Normal"        # drop the quote-prefix style Excel just added
$ws.Range("E41").Value = "  +2.12%  "

# Row 42: MXToken
$ws.Range("E42").Value = "  +1.34%  "

# Row 43: Aave
$ws.Range("D43").Value = "'81.82"   # force text, avoid numeric auto-conversion
$ws.Range("D43").Style = "Normal"        # drop the quote-prefix style Excel just added
$ws.Range("E43").Value = "  -3.13%  "

# Row 44: HuobiToken
$ws.Range("E44").Value = "  +0.34%  "

# Row 45: InjectiveProtocol
$ws.Range("D45").Value = "'13.50"   # force text, avoid numeric auto-conversion
$ws.Range("D45").Style = "Normal"        # drop the quote-prefix style Excel just added
$ws.Range("E45").Value = "  -6.14%  "

# Row 46: FraxShare
$ws.Range("E46").Value = "  -0.48%  "

# Row 47: Kaspa
$ws.Range("E47").Value = "  -4.95%  "

# Row 48: RocketPoolETH
$ws.Range("D48").Value = "1.961.47"
$ws.Range("E48").Value = "  +0.81%  "

# Row 49: Quant
$ws.Range("D49").Value = "'104.65"   # force text, avoid numeric auto-conversion
$ws.Range("D49").Style = "Normal"        # drop the quote-prefix style Excel just added
$ws.Range("E49").Value = "  -0.67%  "

# Row 50: PaxDollar
$ws.Range("E50").Value = "  +0.12%  "

# Row 51: BabyDogeCoin
$ws.Range("D51").Value = "0.0₆0128"
$ws.Range("E51").Value = "  +0.20%  "
